# Fix typo/formatting issues in the "photoStockage" deck.
#  - Slide 7 ("Frontend Architecture"): merge the "Component-based " /
#    "architercture" runs into a single, correctly spelled run
#    "Component-based architecture" (and drop the leftover endParaRPr).
#  - Slide 9 ("Data Management"): collapse the double space in
#    "Tables  of users, ..." down to a single space.

$p = $ppt.ActivePresentation

# --- Slide 7 ("Frontend Architecture"): "Component-based architercture" -> "Component-based architecture"
$s7 = $p.Slides.Item(7)
$tf7 = $s7.Shapes.Item("Content Placeholder 2").TextFrame
$tr7 = $tf7.TextRange

# Remove the whole first paragraph (text + its endParaRPr) outright, then
# insert the corrected, merged text before what is now paragraph 1. This
# mirrors retyping the line from scratch, so the rebuilt run inherits the
# plain run-properties (no stray err="1"/endParaRPr survives the rewrite).
$tr7.Paragraphs(1).Delete()
$null = $tr7.Paragraphs(1).InsertBefore("Component-based architecture" + [char]13)

# --- Slide 9 ("Data Management"): fix the double space before "of users, photos, ..."
$s9 = $p.Slides.Item(9)
$tf9 = $s9.Shapes.Item("Content Placeholder 2").TextFrame
$tr9 = $tf9.TextRange

$dataPara = $tr9.Paragraphs(2)
$firstRun = $tr9.Characters($dataPara.Start, 45)
$firstRun.Text = "Tables of users, photos, likes, categories, "
